# feat: add 2022-Q3 data
#
# The workbook previously had a single "2022-Q2" holdings sheet. We add the
# new "2022-Q3" quarter: the existing "2022-Q2" sheet is renamed to
# "2022-Q3" and repopulated with the Q3 fund-holdings table, while a fresh
# "2022-Q2" sheet is inserted right after it, holding the original Q2 data
# (so nothing is lost, and the sheet keeps its old tab position/relative
# order). The "总计" (totals) summary sheet gets a new Q3 row, with the old
# Q2 totals row pushed down to row 3.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$q2 = $wb.Worksheets.Item(2)

# --- 1. Preserve the existing "2022-Q2" sheet by duplicating it -----------
# The duplicate keeps all of the old fund-holdings rows/formatting intact
# and keeps the "2022-Q2" name; the original sheet object becomes "2022-Q3".
$q2.Copy($null, $q2)
$q2dup = $wb.Worksheets.Item(3)

$q2.Name = "2022-Q3"
$q2dup.Name = "2022-Q2"

# --- 2. Rebuild the (now renamed) "2022-Q3" sheet with the new table ------
$q2.Range("A2:H10").ClearContents()

# Headers (row 1) - copy the bold/centered/bordered header style used on the
# "总计" sheet so the new header row matches it.
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q2.Cells.Item(1, 2 + $i).Value = $headers[$i]
}
$summary.Range("B1").Copy()
$q2.Range("B1:H1").PasteSpecial(-4122)

# Row index column (A2:A10) - same style as the "总计" index column.
$summary.Range("A2").Copy()
$q2.Range("A2:A10").PasteSpecial(-4122)

$data = @(
    @("012526", "广发盛锦混合型证券投资基金A", "25.18", "90.19", "4.13", "1.0399", 6),
    @("010054", "万家健康产业混合A",           "4.60",  "86.96", "3.46", "0.1592", 9),
    @("012527", "广发盛锦混合型证券投资基金C", "1.16",  "90.19", "4.13", "0.0479", 6),
    @("010434", "红土创新医疗保健股票",         "0.44",  "94.49", "6.68", "0.0294", 4),
    @("010055", "万家健康产业混合C",           "0.78",  "86.96", "3.46", "0.0270", 9),
    @("000804", "中信建投稳利混合A",           "0.21",  "38.27", "2.45", "0.0051", 3),
    @("003308", "中信建投睿利灵活配置混合A",   "0.07",  "93.78", "4.39", "0.0031", 4),
    @("006844", "中信建投稳利混合C",           "0.10",  "38.27", "2.45", "0.0024", 3),
    @("004635", "中信建投睿利灵活配置混合C",   "0.03",  "93.78", "4.39", "0.0013", 4)
)

# Columns holding zero-padded fund codes / decimal-text values must stay
# text (otherwise leading zeros are lost and floats drift, e.g. 90.19 ->
# 90.18999999999999), so format them as Text before entering the values.
$q2.Range("B2:B10").NumberFormat = "@"
$q2.Range("D2:G10").NumberFormat = "@"

$row = 2
foreach ($rec in $data) {
    $q2.Cells.Item($row, 1).Value = $row - 2
    $q2.Cells.Item($row, 2).Value = $rec[0]
    $q2.Cells.Item($row, 3).Value = $rec[1]
    $q2.Cells.Item($row, 4).Value = $rec[2]
    $q2.Cells.Item($row, 5).Value = $rec[3]
    $q2.Cells.Item($row, 6).Value = $rec[4]
    $q2.Cells.Item($row, 7).Value = $rec[5]
    $q2.Cells.Item($row, 8).Value = $rec[6]
    $row = $row + 1
}

# --- 3. Update the "总计" summary sheet ------------------------------------
# Push the existing Q2 totals row down to row 3 (values, then formats - see
# note below), then write the new Q3 totals into row 2.
$summary.Range("A2:D2").Copy()
$summary.Range("A3:D3").PasteSpecial(-4163)
$summary.Range("A2:D2").Copy()
$summary.Range("A3:D3").PasteSpecial(-4122)
$summary.Range("A3").Value = 1

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 1.32

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 3
$summary.Range("D3").Value = 0.04

Write-Output "done"
